# Updated results with new RNG run: refresh the 10 trial rows (B2:D11),
# retitle the header/label strings to title-case, drop the stray
# quote-prefixed blank cell at B15, and restore Excel's default page
# margins (file was re-saved from a different machine/Excel build).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / label text -----------------------------------------------
$ws.Range("C1").Value = "Integral"
$ws.Range("D1").Value = "Time"

# --- New trial data (STD / INTEGRAL / TIME columns) ---------------------
$ws.Range("B2").Value = 0.0081905299999999997
$ws.Range("C2").Value = 0.1853822477
$ws.Range("D2").Value = 1.042417639

$ws.Range("B3").Value = 0.0069944886689999996
$ws.Range("C3").Value = 0.18237739
$ws.Range("D3").Value = 1.154465794

$ws.Range("B4").Value = 0.023803789320000001
$ws.Range("C4").Value = 0.20374751260000001
$ws.Range("D4").Value = 1.0620689679999999

$ws.Range("B5").Value = 0.010250486810000001
$ws.Range("C5").Value = 0.1954964075
$ws.Range("D5").Value = 1.006663198

$ws.Range("B6").Value = 0.010444326869999999
$ws.Range("C6").Value = 0.1926358967
$ws.Range("D6").Value = 0.98324969399999995

$ws.Range("B7").Value = 0.0098957974170000001
$ws.Range("C7").Value = 0.2000585584
$ws.Range("D7").Value = 1.000123935

$ws.Range("B8").Value = 0.0071342729589999997
$ws.Range("C8").Value = 0.1909297217
$ws.Range("D8").Value = 0.98166664999999997

$ws.Range("B9").Value = 0.0072722935099999997
$ws.Range("C9").Value = 0.18499853120000001
$ws.Range("D9").Value = 0.98179079000000002

$ws.Range("B10").Value = 0.0059922919610000003
$ws.Range("C10").Value = 0.1840187184
$ws.Range("D10").Value = 0.98115865599999996

$ws.Range("B11").Value = 0.016294530000000002
$ws.Range("C11").Value = 0.1985217933
$ws.Range("D11").Value = 0.98036086200000006

# --- Avg/Std summary rows (formulas stay the same, recalc off new data) -
$ws.Range("A13").Value = "Avg"
$ws.Range("B13").Formula = "=AVERAGE(B2:B11)"
$ws.Range("C13").Formula = "=AVERAGE(C2:C11)"
$ws.Range("D13").Formula = "=AVERAGE(D2:D11)"

$ws.Range("B14").Formula = "=_xlfn.STDEV.S(B2:B11)"
$ws.Range("C14").Formula = "=_xlfn.STDEV.S(C2:C11)"
$ws.Range("D14").Formula = "=_xlfn.STDEV.S(D2:D11)"

# --- Drop the stray quote-prefixed blank cell / row 15 -------------------
$ws.Rows.Item(15).Delete()

# --- Selection cursor moves to where row 15 used to start ---------------
[void]$ws.Range("D15").Select()

# --- Restore Excel's built-in default page margins (inches -> points) ---
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
